# Refresh the "cryptos" price table (GitHub Actions scheduled update).
# Updates Price (D) / Volume(1h) (E) columns with freshly scraped figures,
# and fixes the BNB/Solana rows which had been swapped (rows 5 & 6).
#
# For cells whose new text happens to look like a plain number (e.g.
# "27.15"), the cell's NumberFormat is forced to "@" (Text) right before
# the assignment so Excel's COM layer stores the literal string instead of
# silently re-typing the cell as a float (which would corrupt values like
# "0.0000189" / lose the original text semantics of the source data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Subscript-three character (U+2083) used in the PEPE price string, built
# via [char] so it is not mangled by source encoding.
$sub3 = [string][char]0x2083

$ws.Range("D2").Value = "75.130.62"
$ws.Range("E2").Value = "  +6.01%  "
$ws.Range("D3").Value = "2.703.28"
$ws.Range("E3").Value = "  +9.51%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.61"
$ws.Range("E5").Value = "  +4.43%  "
$ws.Range("B6").Value = "Solana"
$ws.Range("C6").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "186.80"
$ws.Range("E6").Value = "  +9.38%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.539"
$ws.Range("E8").Value = "  +3.71%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.194"
$ws.Range("E9").Value = "  +7.79%  "
$ws.Range("D10").Value = "2.702.24"
$ws.Range("E10").Value = "  +9.61%  "
$ws.Range("E11").Value = "  +0.64%  "
$ws.Range("E12").Value = "  +7.22%  "
$ws.Range("E13").Value = "  +1.37%  "
$ws.Range("D14").Value = "3.202.70"
$ws.Range("E14").Value = "  +9.92%  "
$ws.Range("D15").Value = "75.158.40"
$ws.Range("E15").Value = "  +6.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000189"
$ws.Range("E16").Value = "  +2.65%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.15"
$ws.Range("E17").Value = "  +10.60%  "
$ws.Range("D18").Value = "2.707.78"
$ws.Range("E18").Value = "  +10.52%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.36"
$ws.Range("E19").Value = "  +26.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.09"
$ws.Range("E20").Value = "  +9.89%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "377.72"
$ws.Range("E21").Value = "  +8.35%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.29"
$ws.Range("E22").Value = "  +9.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.11"
$ws.Range("E23").Value = "  +4.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.28"
$ws.Range("E24").Value = "  +3.80%  "
$ws.Range("E25").Value = "  -0.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "70.85"
$ws.Range("E26").Value = "  +6.68%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.19"
$ws.Range("E27").Value = "  +7.40%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.52"
$ws.Range("E28").Value = "  +9.77%  "
$ws.Range("D29").Value = "2.858.91"
$ws.Range("E29").Value = "  +12.26%  "
$ws.Range("E30").Value = "  -0.13%  "
$ws.Range("D31").Value = "0.0" + $sub3 + "0997"
$ws.Range("E31").Value = "  +11.51%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "528.42"
$ws.Range("E32").Value = "  +10.83%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.41"
$ws.Range("E33").Value = "  +10.50%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.85"
$ws.Range("E34").Value = "  +4.50%  "
$ws.Range("E35").Value = "  +9.53%  "
$ws.Range("E36").Value = "  -0.19%  "
$ws.Range("E37").Value = "  +5.50%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "161.71"
$ws.Range("E38").Value = "  +2.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.51"
$ws.Range("E39").Value = "  +5.79%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19.32"
$ws.Range("E40").Value = "  +0.77%  "
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "174.88"
$ws.Range("E42").Value = "  +25.54%  "
$ws.Range("E43").Value = "  +11.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.70"
$ws.Range("E44").Value = "  +8.71%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.332"
$ws.Range("E45").Value = "  +7.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.24"
$ws.Range("E46").Value = "  +11.93%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.42"
$ws.Range("E47").Value = "  +10.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "39.26"
$ws.Range("E48").Value = "  +2.76%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0848"
$ws.Range("E49").Value = "  +16.53%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.70"
$ws.Range("E50").Value = "  +7.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.547"
$ws.Range("E51").Value = "  +9.62%  "
